# Generate Report for Handoff
# Regenerates handoff artifacts for the localization status report:
# a new source file UUID replaces the old one, fresh handoff xliff
# files/timestamps are recorded, and the (not-yet-produced) handback
# file/target file/timestamp fields are reset to their "not done yet"
# defaults.

$oldGuid = "f6bf5b5f-04f3-493a-8e89-5ae441e222d8"
$newGuid = "54cf65f1-86e3-4a25-be0e-683f11feb21e"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = "2016-09-01 17:10:10"

foreach ($hl in $ws1.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$2") {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$ws2.Range("A2").Value = "$newGuid.md"
foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

$ws2.Range("G2").Value = "$newGuid.30ebfddd2b2acba1497ddfbed9a04f55dfd12bc4.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-01 17:09:57"

# Latest Target File (I2) no longer has a handback, so its hyperlink is
# removed and the cell is cleared.
foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$I`$2") {
        $hl.Delete()
    }
}
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"

# Latest Handback File (J2) cleared - no handback produced yet.
$ws2.Range("J2").Value = ""

# Latest Handback DateTime (K2) reset to the "never happened" sentinel.
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Columns.Item(9).ColumnWidth = 17.83
$ws2.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$ws3.Range("A2").Value = "$newGuid.md"
foreach ($hl in $ws3.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

$ws3.Range("G2").Value = "$newGuid.30ebfddd2b2acba1497ddfbed9a04f55dfd12bc4.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-01 17:10:10"

# Latest Target File (I2) no longer has a handback, so its hyperlink is
# removed and the cell is cleared.
foreach ($hl in $ws3.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$I`$2") {
        $hl.Delete()
    }
}
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"

# Latest Handback File (J2) cleared - no handback produced yet.
$ws3.Range("J2").Value = ""

# Latest Handback DateTime (K2) reset to the "never happened" sentinel.
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Columns.Item(9).ColumnWidth = 17.83
$ws3.Columns.Item(10).ColumnWidth = 20.83
